$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1) and "全部类型" sheet (sheet4): F4 75 -> 76, F5 294 -> 295
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 76
    $ws.Range("F5").Value = 295
}
